$d = $word.ActiveDocument

# Locate the final paragraph (ilvl 0, numId 1) that currently holds just a
# placeholder space and the "_GoBack" bookmark.
$lastPara = $d.Paragraphs.Last
$lastRange = $lastPara.Range

# Drop the existing "_GoBack" bookmark; it will be re-created further down,
# anchored to the new "...相关定义" paragraph, matching the target layout.
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# Re-fetch the paragraph/range now that the bookmark is gone.
$lastPara = $d.Paragraphs.Last
$lastRange = $lastPara.Range
$contentRange = $d.Range($lastRange.Start, $lastRange.End - 1)

# Replace the lone placeholder run (" ") with the "理论基础" heading text,
# carrying the eastAsia font hint exactly like the rest of the document.
$run1Xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>理论基础</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$contentRange.InsertXML($run1Xml)

# Append the three new paragraphs (reference, radiative-transfer bullet, and
# the trailing empty bullet) right after the paragraph we just edited.
$lastPara = $d.Paragraphs.Last
$lastRange = $lastPara.Range
$insertPoint = $d.Range($lastRange.End, $lastRange.End)
$tailXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="a3"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:ind w:firstLineChars="0"/></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>Cao</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>，</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>20</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>19</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="a3"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr><w:ind w:firstLineChars="0"/></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>辐射传输方程</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>、</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>LST</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>相关定义</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p><w:pPr><w:pStyle w:val="a3"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr><w:ind w:firstLineChars="0"/><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertPoint.InsertXML($tailXml)

Write-Output "edit complete"
